# Add Steve's user stories for the Administrator role.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Update Map Data -------------------------------------------------
$ws.Range("C9").Value = "Administrator"
$ws.Range("D9").Value = "Designate points on map that link with stories."
$ws.Range("E9").Value = "To provide users an easy way to view the location to which a story is connected."
$ws.Range("F9").Value = "A story has been posted."
$ws.Range("G9").Value = "A point on the map is linked with the story."
$ws.Range("C9:G9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 60

# --- Row 10: Delete or Edit Post/Comment ------------------------------------
$ws.Range("C10").Value = "Administrator"
$ws.Range("D10").Value = "Delete or edit submissions by users and content publishers."
$ws.Range("E10").Value = "To provide the ability to moderate posts and comments in case inappropriate content is submitted."
$ws.Range("F10").Value = "A post/comment has been submitted and an administrator has a reason to edit or delete it."
$ws.Range("G10").Value = "The post/comment has been changed or deleted."
$ws.Range("C10:G10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 90

# --- Row 11: View Analytics --------------------------------------------------
$ws.Range("C11").Value = "Administrator"
$ws.Range("D11").Value = "View analytics about the website such as demographics and page views"
$ws.Range("E11").Value = "To give administrators meaningful statistics about the website."
$ws.Range("F11").Value = "Analytics have been collected and administrator is logged in."
$ws.Range("G11").Value = "The administrator is able to view analytics."
$ws.Range("C11:G11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 120

# --- Update the active selection to match the authored workbook -------------
$ws.Range("E6").Select()
